$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.082.17"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.088.00"
$ws.Range("E3").Value = "  +9.00%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.67"
$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.660"
$ws.Range("E6").Value = "  -4.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.40"
$ws.Range("E8").Value = "  +7.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.58"
$ws.Range("E9").Value = "  +4.31%  "

$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0745"
$ws.Range("E11").Value = "  -2.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.107"
$ws.Range("E12").Value = "  +6.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.57"
$ws.Range("E13").Value = "  -4.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.390.80"
$ws.Range("E14").Value = "  +9.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.833"
$ws.Range("E15").Value = "  +1.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.091.89"
$ws.Range("E16").Value = "  +9.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.15"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.997.53"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("E19").Value = "  -2.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0826"
$ws.Range("E20").Value = "  -3.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.35"
$ws.Range("E21").Value = "  -2.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.04"
$ws.Range("E22").Value = "  -4.12%  "

$ws.Range("E23").Value = "  +2.06%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("E25").Value = "  +2.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.83"
$ws.Range("E26").Value = "  +1.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.35"
$ws.Range("E27").Value = "  +5.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.07"
$ws.Range("E28").Value = "  +12.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("E29").Value = "  -8.38%  "

$ws.Range("B30").Value = "Gas"
$ws.Range("C30").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.23"
$ws.Range("E30").Value = "  +35.52%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.123"
$ws.Range("E31").Value = "  -4.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  +28.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.51"
$ws.Range("E33").Value = "  -1.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0611"
$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0929"
$ws.Range("E35").Value = "  +3.22%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.32"
$ws.Range("E37").Value = "  +16.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.12"
$ws.Range("E38").Value = "  -4.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.80"
$ws.Range("E39").Value = "  -5.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.33"
$ws.Range("E40").Value = "  -10.25%  "

$ws.Range("E41").Value = "  -1.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.68"
$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.16"
$ws.Range("E43").Value = "  +6.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.68"
$ws.Range("E44").Value = "  -6.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("E45").Value = "  -2.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.349.20"
$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0865"
$ws.Range("E47").Value = "  +4.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.94"
$ws.Range("E48").Value = "  +4.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.04"
$ws.Range("E49").Value = "  +9.86%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.276.87"
$ws.Range("E50").Value = "  +8.71%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.28"
$ws.Range("E51").Value = "  -4.96%  "

